$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Remove old password-manager rows, keep header row 1 ---
$ws.Range("A4:R16").ClearContents()

# --- Row 1: bump counter value, keep header texts as-is ---
$ws.Range("F1").Value = 4

# --- Row 2: new PrivateKey label + value ---
$ws.Range("E2").Value = "PrivateKey:"
$ws.Range("F2").Value = "57a0!c0>u2.r4(c9-r4#l2^x0~O7*P3}K9-C9#m8,w4+z0*Y7+A4-P8*I7{v1^=114a0!107c0>138u2.111r4(119c9-179r4#110l2^100x0~177O7*107P3}150K9-128C9#113m8,183w4+162z0*105Y7+181A4-128P8*110I7{102v1^=191Q6>189z5+201x2]197G8~200i0|203G5-="

# --- Row 3: new "camilo" entry ---
$ws.Range("A3").Value = "camilo"
$ws.Range("B3").Value = "72b4]t9.F8}O6,N9]Y8{K0|Y3#x1]k3|s8)T2!P6:N1&O8(c4!C5|n3,v1`i2]^128b4]128t9.147F8}120O6,110N9]192Y8{128K0|197Y3#137x1]127k3|109s8)158T2!125P6:133N1&184O8(129c4!163C5|160n3,126v1`196i2]^187K0&185x8}197A7%193V5=196p9[199X6|^"
$ws.Range("C3").Value = "20w5|c6/Q2]D0_G9)x6}H2!Z4#K2_e8<I6-h7+P7.i4*g8}k9.u3,B4[C5)D0]@71w5|70c6/107Q2]70D0_66G9)123x6}69H2!61Z4#129K2_71e8<78I6-124h7+69P7.78i4*90g8}69k9.81u3,100B4[75C5)61D0]@131W2.129g1)141m3:137h1:140F1=143P7)@"

# --- New column formatting for column E ---
$ws.Columns("E").ColumnWidth = 12.5703125

# --- Row 11: empty cell with underline font style ---
$ws.Range("F11").Font.Underline = -4119

# --- Update selection to F1 ---
$ws.Range("F1").Select()
